# Daily attendance processing - 2025-10-30 23:42:17
# Applies updated attendance figures / recorded-session info to the
# "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Helper: write a text value into a cell without Excel "smart"
# auto-converting percent-looking / numeric-looking strings into
# numbers. We force a text number format, assign the value, then
# restore the original cell formatting by pasting formats back in
# from an UNTOUCHED reference cell that already carries the
# desired (original) style.
# ---------------------------------------------------------------
function Set-TextValue($cell, $value, $formatSourceCell) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $formatSourceCell.Copy()
    $cell.PasteSpecial(-4122)
}

# =================================================================
# Class / overall statistics (K:L columns)
# =================================================================
$ws.Range("L6").Value = 194
$ws.Range("L7").Value = 9

Set-TextValue $ws.Range("L9") "72.9%" $ws.Range("K9")
Set-TextValue $ws.Range("L10") "71.7%" $ws.Range("K10")

# =================================================================
# Group statistics (row 16 -> B1B, row 19 -> B1E)
# =================================================================
$ws.Range("O16").Value = 42
$ws.Range("P16").Value = 2
Set-TextValue $ws.Range("R16") "79.2%" $ws.Range("N16")
Set-TextValue $ws.Range("S16") "71.8%" $ws.Range("N16")

$ws.Range("O19").Value = 40
$ws.Range("P19").Value = 3
Set-TextValue $ws.Range("R19") "75.5%" $ws.Range("N19")
Set-TextValue $ws.Range("S19") "67.3%" $ws.Range("N19")

# =================================================================
# Sessions that moved from "Not Recorded" (pink) to "Recorded"
# (normal) status. Copy the look/format of an already-recorded
# row in the same block, then fill in the actual attendance data.
# =================================================================

# --- B1B / CARDIOLOGY block (students out of 65), template row 57 ---
$ws.Range("A57:I57").Copy()
$ws.Range("A56:I56").PasteSpecial(-4122)
$ws.Range("G56").Value = "rana.s.nasser1999@gmail.com"
$ws.Range("H56").Value = "65/65"
$ws.Range("I56").Value = "Recorded"

$ws.Range("A57:I57").Copy()
$ws.Range("A58:I58").PasteSpecial(-4122)
$ws.Range("G58").Value = "rana.s.nasser1999@gmail.com"
$ws.Range("H58").Value = "60/65"
$ws.Range("I58").Value = "Recorded"

$ws.Range("A57:I57").Copy()
$ws.Range("A60:I60").PasteSpecial(-4122)
$ws.Range("G60").Value = "rana.s.nasser1999@gmail.com"
$ws.Range("H60").Value = "60/65"
$ws.Range("I60").Value = "Recorded"

$ws.Range("A57:I57").Copy()
$ws.Range("A61:I61").PasteSpecial(-4122)
$ws.Range("G61").Value = "rana.s.nasser1999@gmail.com"
$ws.Range("H61").Value = "65/65"
$ws.Range("I61").Value = "Recorded"

$ws.Range("A57:I57").Copy()
$ws.Range("A63:I63").PasteSpecial(-4122)
$ws.Range("G63").Value = "rana.s.nasser1999@gmail.com"
$ws.Range("H63").Value = "65/65"
$ws.Range("I63").Value = "Recorded"

$ws.Range("A57:I57").Copy()
$ws.Range("A64:I64").PasteSpecial(-4122)
$ws.Range("G64").Value = "rana.s.nasser1999@gmail.com"
$ws.Range("H64").Value = "65/65"
$ws.Range("I64").Value = "Recorded"

# --- B1E / CARDIOLOGY block (students out of 70), template row 217 ---
$ws.Range("A217:I217").Copy()
$ws.Range("A218:I218").PasteSpecial(-4122)
$ws.Range("G218").Value = "rana.s.nasser1999@gmail.com"
$ws.Range("H218").Value = "70/70"
$ws.Range("I218").Value = "Recorded"

$ws.Range("A217:I217").Copy()
$ws.Range("A219:I219").PasteSpecial(-4122)
$ws.Range("G219").Value = "rana.s.nasser1999@gmail.com"
$ws.Range("H219").Value = "70/70"
$ws.Range("I219").Value = "Recorded"

$ws.Range("A217:I217").Copy()
$ws.Range("A223:I223").PasteSpecial(-4122)
$ws.Range("G223").Value = "rana.s.nasser1999@gmail.com"
$ws.Range("H223").Value = "70/70"
$ws.Range("I223").Value = "Recorded"

# --- B1E / IMMUNO&HEMA block (students out of 70), template row 246 ---
$ws.Range("A246:I246").Copy()
$ws.Range("A245:I245").PasteSpecial(-4122)
$ws.Range("G245").Value = "System"
$ws.Range("H245").Value = "30/70"
$ws.Range("I245").Value = "Recorded"

# =================================================================
# "Recorded By" (column G) entries whose listed contributors were
# reordered.
# =================================================================
$ws.Range("G83").Value = "System, Mai.elsebaie@gmail.com"

$ws.Range("G85").Value = "emp17066@med.asu.edu.eg, 160392@med.asu.edu.com"
$ws.Range("G88").Value = "emp17066@med.asu.edu.eg, 160392@med.asu.edu.com"
$ws.Range("G91").Value = "emp17066@med.asu.edu.eg, 160392@med.asu.edu.com"
$ws.Range("G92").Value = "emp17066@med.asu.edu.eg, 160392@med.asu.edu.com"
$ws.Range("G248").Value = "emp17066@med.asu.edu.eg, 160392@med.asu.edu.com"

$ws.Range("G120").Value = "System, ahmedali78112@gmail.com"
$ws.Range("G121").Value = "System, ahmedali78112@gmail.com"
$ws.Range("G122").Value = "System, ahmedali78112@gmail.com"
$ws.Range("G123").Value = "System, ahmedali78112@gmail.com"
$ws.Range("G124").Value = "System, ahmedali78112@gmail.com"
$ws.Range("G125").Value = "System, ahmedali78112@gmail.com"
$ws.Range("G126").Value = "System, ahmedali78112@gmail.com"
$ws.Range("G127").Value = "System, ahmedali78112@gmail.com"
